$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.158.47'
$ws.Range('E2').Value = '  -0.26%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.048.39'
$ws.Range('E3').Value = '  -0.73%  '
$ws.Range('E4').Value = '  -0.42%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.94'
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.21'
$ws.Range('E6').Value = '  -0.41%  '
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.538'
$ws.Range('E8').Value = '  -1.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.046.07'
$ws.Range('E9').Value = '  -0.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.153'
$ws.Range('E10').Value = '  -0.85%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.79'
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('E12').Value = '  -1.86%  '
$ws.Range('E13').Value = '  -1.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.28'
$ws.Range('E14').Value = '  -1.67%  '
$ws.Range('E15').Value = '  +1.75%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.552.49'
$ws.Range('E16').Value = '  -0.75%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.15'
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.178.33'
$ws.Range('E18').Value = '  -0.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.053.10'
$ws.Range('E19').Value = '  -0.67%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '476.84'
$ws.Range('E20').Value = '  +1.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.26'
$ws.Range('E21').Value = '  -2.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.705'
$ws.Range('E22').Value = '  -1.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.51'
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('E24').Value = '  +2.37%  '
$ws.Range('E25').Value = '  +1.34%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.70'
$ws.Range('E26').Value = '  -2.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.67'
$ws.Range('E27').Value = '  +7.74%  '
$ws.Range('E28').Value = '  +0.35%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.35'
$ws.Range('E29').Value = '  +0.96%  '
$ws.Range('E30').Value = '  +0.16%  '
$ws.Range('E31').Value = '  -0.34%  '
$ws.Range('E32').Value = '  +0.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.67'
$ws.Range('E33').Value = '  +2.00%  '
$ws.Range('E34').Value = '  -2.30%  '
$ws.Range('E35').Value = '  +1.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0₃0817'
$ws.Range('E36').Value = '  -3.37%  '
$ws.Range('E37').Value = '  +1.16%  '
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.24'
$ws.Range('E38').Value = '  -2.68%  '
$ws.Range('B39').Value = 'Filecoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.89'
$ws.Range('E39').Value = '  -2.91%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.25'
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('E41').Value = '  +0.48%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '434.71'
$ws.Range('E42').Value = '  -2.18%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.287'
$ws.Range('E43').Value = '  +1.41%  '
$ws.Range('E44').Value = '  +3.11%  '
$ws.Range('E45').Value = '  +0.47%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.827.92'
$ws.Range('E46').Value = '  +1.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '38.29'
$ws.Range('E47').Value = '  -4.48%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '129.84'
$ws.Range('E48').Value = '  -0.88%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '25.08'
$ws.Range('E50').Value = '  +0.18%  '
$ws.Range('E51').Value = '  -1.31%  '
